$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7
$ws.Cells.Item(7, 11).Value = 9
$ws.Cells.Item(7, 15).Value = 1.7
$ws.Cells.Item(7, 18).Value = 1.91
$ws.Cells.Item(7, 19).Value = 1.91

# Row 8
$ws.Cells.Item(8, 14).Value = 1.73
$ws.Cells.Item(8, 15).Value = 2.1
$ws.Cells.Item(8, 19).Value = 1.7

# Row 9
$ws.Cells.Item(9, 14).Value = 1.67
$ws.Cells.Item(9, 15).Value = 2.2
$ws.Cells.Item(9, 18).Value = 1.57

# Row 10
$ws.Cells.Item(10, 14).Value = 2.63
$ws.Cells.Item(10, 15).Value = 1.5
$ws.Cells.Item(10, 19).Value = 1.67

# Row 11
$ws.Cells.Item(11, 18).Value = 1.91
$ws.Cells.Item(11, 19).Value = 1.91
$ws.Cells.Item(11, 28).Value = 17

# Row 15
$ws.Cells.Item(15, 7).Value = 1.85
$ws.Cells.Item(15, 8).Value = 3.25
$ws.Cells.Item(15, 10).Value = 1.05
$ws.Cells.Item(15, 11).Value = 11
$ws.Cells.Item(15, 12).Value = 1.25
$ws.Cells.Item(15, 14).Value = 1.83
$ws.Cells.Item(15, 15).Value = 2.03
$ws.Cells.Item(15, 16).Value = 1.36
$ws.Cells.Item(15, 17).Value = 3
$ws.Cells.Item(15, 19).Value = 2.05
$ws.Cells.Item(15, 20).Value = 8
$ws.Cells.Item(15, 21).Value = 9.5
$ws.Cells.Item(15, 22).Value = 8.5
$ws.Cells.Item(15, 23).Value = 17
$ws.Cells.Item(15, 24).Value = 15
$ws.Cells.Item(15, 25).Value = 23
$ws.Cells.Item(15, 26).Value = 11
$ws.Cells.Item(15, 27).Value = 6.5
$ws.Cells.Item(15, 28).Value = 13
$ws.Cells.Item(15, 29).Value = 41
$ws.Cells.Item(15, 30).Value = 151
$ws.Cells.Item(15, 31).Value = 13
$ws.Cells.Item(15, 32).Value = 21
$ws.Cells.Item(15, 33).Value = 13
$ws.Cells.Item(15, 34).Value = 41
$ws.Cells.Item(15, 35).Value = 34
$ws.Cells.Item(15, 36).Value = 34

# Row 16
$ws.Cells.Item(16, 7).Value = 3.2
$ws.Cells.Item(16, 8).Value = 3.5
$ws.Cells.Item(16, 9).Value = 2
$ws.Cells.Item(16, 10).Value = 1.03
$ws.Cells.Item(16, 11).Value = 15
$ws.Cells.Item(16, 12).Value = 1.2
$ws.Cells.Item(16, 13).Value = 4.33
$ws.Cells.Item(16, 14).Value = 1.67
$ws.Cells.Item(16, 15).Value = 2.15
$ws.Cells.Item(16, 16).Value = 1.33
$ws.Cells.Item(16, 17).Value = 3.25
$ws.Cells.Item(16, 18).Value = 1.57
$ws.Cells.Item(16, 19).Value = 2.25
$ws.Cells.Item(16, 20).Value = 12
$ws.Cells.Item(16, 21).Value = 17
$ws.Cells.Item(16, 22).Value = 11
$ws.Cells.Item(16, 23).Value = 34
$ws.Cells.Item(16, 24).Value = 23
$ws.Cells.Item(16, 25).Value = 26
$ws.Cells.Item(16, 26).Value = 15
$ws.Cells.Item(16, 27).Value = 7
$ws.Cells.Item(16, 28).Value = 12
$ws.Cells.Item(16, 29).Value = 41
$ws.Cells.Item(16, 30).Value = 126
$ws.Cells.Item(16, 31).Value = 10
$ws.Cells.Item(16, 32).Value = 12
$ws.Cells.Item(16, 34).Value = 21
$ws.Cells.Item(16, 35).Value = 15
$ws.Cells.Item(16, 36).Value = 21

# Row 17
$ws.Cells.Item(17, 7).Value = 2.82
$ws.Cells.Item(17, 8).Value = 3
$ws.Cells.Item(17, 9).Value = 2.47
$ws.Cells.Item(17, 18).Value = 1.83
$ws.Cells.Item(17, 19).Value = 1.78
$ws.Cells.Item(17, 20).Value = 7.5
$ws.Cells.Item(17, 26).Value = 7.6
$ws.Cells.Item(17, 27).Value = 5.9
$ws.Cells.Item(17, 28).Value = 15
$ws.Cells.Item(17, 29).Value = 80
$ws.Cells.Item(17, 31).Value = 7.2
$ws.Cells.Item(17, 32).Value = 11.75
$ws.Cells.Item(17, 34).Value = 27

# Row 21
$ws.Cells.Item(21, 7).Value = 7
$ws.Cells.Item(21, 9).Value = 1.42
$ws.Cells.Item(21, 18).Value = 2.1
$ws.Cells.Item(21, 19).Value = 1.67
$ws.Cells.Item(21, 20).Value = 17
$ws.Cells.Item(21, 22).Value = 21
$ws.Cells.Item(21, 23).Value = 81
$ws.Cells.Item(21, 27).Value = 8.5
$ws.Cells.Item(21, 35).Value = 12

# Row 23
$ws.Cells.Item(23, 14).Value = 2.08
$ws.Cells.Item(23, 15).Value = 1.73

# Row 24
$ws.Cells.Item(24, 10).Value = 1.09
$ws.Cells.Item(24, 11).Value = 6
$ws.Cells.Item(24, 30).Value = 101

# Row 45
$ws.Cells.Item(45, 33).Value = 11.75
$ws.Cells.Item(45, 34).Value = 45

# Row 46
$ws.Cells.Item(46, 13).Value = 2.77
$ws.Cells.Item(46, 14).Value = 2.15
$ws.Cells.Item(46, 15).Value = 1.62
$ws.Cells.Item(46, 18).Value = 1.93
$ws.Cells.Item(46, 19).Value = 1.78
$ws.Cells.Item(46, 22).Value = 13
$ws.Cells.Item(46, 25).Value = 45
$ws.Cells.Item(46, 29).Value = 90
$ws.Cells.Item(46, 31).Value = 6
$ws.Cells.Item(46, 32).Value = 8.25
$ws.Cells.Item(46, 33).Value = 8.5
$ws.Cells.Item(46, 35).Value = 17.5

# Row 49
$ws.Cells.Item(49, 20).Value = 7.1
$ws.Cells.Item(49, 23).Value = 20
$ws.Cells.Item(49, 29).Value = 70
$ws.Cells.Item(49, 30).Value = 600
$ws.Cells.Item(49, 34).Value = 45

# Row 52
$ws.Cells.Item(52, 7).Value = 2.67
$ws.Cells.Item(52, 8).Value = 3.4
$ws.Cells.Item(52, 9).Value = 2.37
$ws.Cells.Item(52, 16).Value = 1.39
$ws.Cells.Item(52, 17).Value = 2.55
$ws.Cells.Item(52, 18).Value = 1.72
$ws.Cells.Item(52, 19).Value = 1.88
$ws.Cells.Item(52, 20).Value = 8.75
$ws.Cells.Item(52, 21).Value = 13.5
$ws.Cells.Item(52, 22).Value = 10
$ws.Cells.Item(52, 23).Value = 29
$ws.Cells.Item(52, 24).Value = 23
$ws.Cells.Item(52, 26).Value = 10
$ws.Cells.Item(52, 27).Value = 6.6
$ws.Cells.Item(52, 28).Value = 14.5
$ws.Cells.Item(52, 29).Value = 70
$ws.Cells.Item(52, 31).Value = 8
$ws.Cells.Item(52, 32).Value = 11.5
$ws.Cells.Item(52, 33).Value = 9.5
$ws.Cells.Item(52, 34).Value = 24
$ws.Cells.Item(52, 35).Value = 19.5
$ws.Cells.Item(52, 36).Value = 30

# Row 63
$ws.Cells.Item(63, 7).Value = 2.6
$ws.Cells.Item(63, 8).Value = 3.2
$ws.Cells.Item(63, 10).Value = 1.05
$ws.Cells.Item(63, 11).Value = 11
$ws.Cells.Item(63, 12).Value = 1.25
$ws.Cells.Item(63, 13).Value = 3.75
$ws.Cells.Item(63, 14).Value = 1.85
$ws.Cells.Item(63, 15).Value = 2
$ws.Cells.Item(63, 16).Value = 1.36
$ws.Cells.Item(63, 17).Value = 3
$ws.Cells.Item(63, 18).Value = 1.67
$ws.Cells.Item(63, 20).Value = 10
$ws.Cells.Item(63, 21).Value = 15
$ws.Cells.Item(63, 23).Value = 26
$ws.Cells.Item(63, 24).Value = 21
$ws.Cells.Item(63, 25).Value = 29
$ws.Cells.Item(63, 26).Value = 11
$ws.Cells.Item(63, 29).Value = 41
$ws.Cells.Item(63, 30).Value = 151
$ws.Cells.Item(63, 31).Value = 9.5
$ws.Cells.Item(63, 32).Value = 13
$ws.Cells.Item(63, 33).Value = 10
$ws.Cells.Item(63, 34).Value = 23
$ws.Cells.Item(63, 36).Value = 26

# Row 65
$ws.Cells.Item(65, 7).Value = 2.35
$ws.Cells.Item(65, 9).Value = 2.9
$ws.Cells.Item(65, 10).Value = 1.02
$ws.Cells.Item(65, 11).Value = 11
$ws.Cells.Item(65, 22).Value = 9.5
$ws.Cells.Item(65, 23).Value = 21
$ws.Cells.Item(65, 35).Value = 23

# Row 66
$ws.Cells.Item(66, 7).Value = 5.25
$ws.Cells.Item(66, 8).Value = 4
$ws.Cells.Item(66, 9).Value = 1.6
$ws.Cells.Item(66, 10).Value = 1.01
$ws.Cells.Item(66, 11).Value = 13
$ws.Cells.Item(66, 29).Value = 41
$ws.Cells.Item(66, 32).Value = 8.5

# Row 68
$ws.Cells.Item(68, 7).Value = 4.2
$ws.Cells.Item(68, 8).Value = 3.5
$ws.Cells.Item(68, 13).Value = 3.6
$ws.Cells.Item(68, 19).Value = 2.05
$ws.Cells.Item(68, 20).Value = 13.5
$ws.Cells.Item(68, 21).Value = 26
$ws.Cells.Item(68, 25).Value = 37
$ws.Cells.Item(68, 28).Value = 13.5
$ws.Cells.Item(68, 29).Value = 55
$ws.Cells.Item(68, 31).Value = 7.6

# Row 69
$ws.Cells.Item(69, 8).Value = 4.1
$ws.Cells.Item(69, 9).Value = 5.8
$ws.Cells.Item(69, 20).Value = 7.2
$ws.Cells.Item(69, 21).Value = 7.1
$ws.Cells.Item(69, 24).Value = 11.5
$ws.Cells.Item(69, 35).Value = 65

# Row 70
$ws.Cells.Item(70, 9).Value = 1.39
$ws.Cells.Item(70, 16).Value = 1.28
$ws.Cells.Item(70, 22).Value = 20

# Row 71
$ws.Cells.Item(71, 7).Value = 1.88
$ws.Cells.Item(71, 8).Value = 3.85
$ws.Cells.Item(71, 10).Value = 1.04
$ws.Cells.Item(71, 11).Value = 9
$ws.Cells.Item(71, 12).Value = 1.22
$ws.Cells.Item(71, 13).Value = 3.95
$ws.Cells.Item(71, 14).Value = 1.65
$ws.Cells.Item(71, 15).Value = 2.15
$ws.Cells.Item(71, 16).Value = 1.34
$ws.Cells.Item(71, 17).Value = 3.1
$ws.Cells.Item(71, 18).Value = 1.62
$ws.Cells.Item(71, 19).Value = 2.15
$ws.Cells.Item(71, 20).Value = 8.5
$ws.Cells.Item(71, 21).Value = 10.5
$ws.Cells.Item(71, 24).Value = 14.5
$ws.Cells.Item(71, 25).Value = 24
$ws.Cells.Item(71, 26).Value = 9
$ws.Cells.Item(71, 27).Value = 7.9
$ws.Cells.Item(71, 28).Value = 14.5
$ws.Cells.Item(71, 29).Value = 60
$ws.Cells.Item(71, 30).Value = 400
$ws.Cells.Item(71, 31).Value = 12
$ws.Cells.Item(71, 36).Value = 35

# Row 73
$ws.Cells.Item(73, 10).Value = 1.03
$ws.Cells.Item(73, 11).Value = 15

# Row 75
$ws.Cells.Item(75, 12).Value = 1.36
$ws.Cells.Item(75, 13).Value = 3

# Row 77
$ws.Cells.Item(77, 7).Value = 8.75
$ws.Cells.Item(77, 8).Value = 4.45
$ws.Cells.Item(77, 9).Value = 1.35
$ws.Cells.Item(77, 11).Value = 8.5
$ws.Cells.Item(77, 12).Value = 1.2
$ws.Cells.Item(77, 13).Value = 4.1
$ws.Cells.Item(77, 14).Value = 1.6
$ws.Cells.Item(77, 15).Value = 2.2
$ws.Cells.Item(77, 16).Value = 1.33
$ws.Cells.Item(77, 17).Value = 3.05
$ws.Cells.Item(77, 18).Value = 1.83
$ws.Cells.Item(77, 19).Value = 1.87
$ws.Cells.Item(77, 20).Value = 23
$ws.Cells.Item(77, 21).Value = 65
$ws.Cells.Item(77, 22).Value = 26
$ws.Cells.Item(77, 23).Value = 250
$ws.Cells.Item(77, 24).Value = 100
$ws.Cells.Item(77, 25).Value = 75
$ws.Cells.Item(77, 26).Value = 8.5
$ws.Cells.Item(77, 27).Value = 9
$ws.Cells.Item(77, 28).Value = 18.5
$ws.Cells.Item(77, 29).Value = 80
$ws.Cells.Item(77, 30).Value = 600
$ws.Cells.Item(77, 31).Value = 7.6
$ws.Cells.Item(77, 32).Value = 6.9
$ws.Cells.Item(77, 33).Value = 8
$ws.Cells.Item(77, 34).Value = 9
$ws.Cells.Item(77, 35).Value = 10.5

# Row 80
$ws.Cells.Item(80, 20).Value = 7.7
$ws.Cells.Item(80, 21).Value = 10
$ws.Cells.Item(80, 31).Value = 10.75
$ws.Cells.Item(80, 34).Value = 45

# Row 81
$ws.Cells.Item(81, 20).Value = 6.9
$ws.Cells.Item(81, 24).Value = 15.5
$ws.Cells.Item(81, 30).Value = 700
$ws.Cells.Item(81, 31).Value = 10
$ws.Cells.Item(81, 32).Value = 18.5
$ws.Cells.Item(81, 36).Value = 45
